# Generate Report for Handoff
# Update the localization status report:
#  - Mark rows 8,9,10,12,13,14 as priority "ht" (handoff type) on the
#    zh-cn and de-de sheets.
#  - Bump the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#    timestamps for those same rows to reflect the newly generated
#    handoff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: Latest HO Xliff Generate Date (column G)
    $overview.Range("G$r").Value = "2016-09-05 16:25:40"

    # zh-cn sheet: Priority (column E) + Latest Handoff Datetime (column H)
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-05 16:25:36"

    # de-de sheet: Priority (column E) + Latest Handoff Datetime (column H)
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-05 16:25:40"
}
